$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2..201).
# The sheet was refreshed one day later, so bump each value by one day
# (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04).
$lastRow = 201
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current + 1
    }
}
